# The deck's single slide master (theme/theme2.xml, "Integral") is swapped
# with the unused notes-master theme (theme/theme1.xml, "Office Theme") in
# the target revision: the slide master ends up showing the stock "Office"
# colour palette instead of "Integral".
#
# The color values below are the twelve DrawingML theme colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) of the stock "Office"
# theme, expressed as COM RGB() packed integers (0x00BBGGRR, i.e.
# r + g*256 + b*65536) because PowerPoint's ThemeColor.RGB property takes/
# returns colors in that packed form.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    @{ Index = 1;  Name = "dk1";      R = 0x00; G = 0x00; B = 0x00 },
    @{ Index = 2;  Name = "lt1";      R = 0xFF; G = 0xFF; B = 0xFF },
    @{ Index = 3;  Name = "dk2";      R = 0x44; G = 0x54; B = 0x6A },
    @{ Index = 4;  Name = "lt2";      R = 0xE7; G = 0xE6; B = 0xE6 },
    @{ Index = 5;  Name = "accent1";  R = 0x5B; G = 0x9B; B = 0xD5 },
    @{ Index = 6;  Name = "accent2";  R = 0xED; G = 0x7D; B = 0x31 },
    @{ Index = 7;  Name = "accent3";  R = 0xA5; G = 0xA5; B = 0xA5 },
    @{ Index = 8;  Name = "accent4";  R = 0xFF; G = 0xC0; B = 0x00 },
    @{ Index = 9;  Name = "accent5";  R = 0x44; G = 0x72; B = 0xC4 },
    @{ Index = 10; Name = "accent6";  R = 0x70; G = 0xAD; B = 0x47 },
    @{ Index = 11; Name = "hlink";    R = 0x05; G = 0x63; B = 0xC1 },
    @{ Index = 12; Name = "folHlink"; R = 0x95; G = 0x4F; B = 0x72 }
)

foreach ($entry in $officeColors) {
    $packed = $entry.R + ($entry.G * 256) + ($entry.B * 65536)
    $themeColor = $tcs.Colors($entry.Index)
    $themeColor.RGB = $packed
}
